$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (displayed name); underlying tab title maps to sheet name in OOXML.
$ws.Name = "GammaFiber2F"

# Fix tiny rounding differences in row 13
$ws.Range("C13").Value = 0.9937378308183796
$ws.Range("F13").Value = 0.9937378308183796

# Add new row 16 with HexGrid-60degTilt5degRes data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.045273451079792
$ws.Range("D16").Value = 0.8527933969363093
$ws.Range("E16").Value = 1.019274105624675
$ws.Range("F16").Value = 1.045273451079792
$ws.Range("G16").Value = 0.920361292935602
$ws.Range("H16").Value = 1.065469375697556
$ws.Range("I16").Value = 1.02580319011694
$ws.Range("J16").Value = 0.8527933969363093
$ws.Range("K16").Value = 0.9360337512804924
$ws.Range("L16").Value = 0.9906536011801422
$ws.Range("M16").Value = 0.9881624687318125

# Match style of A column cells (border/bold/centered) used on A2:A15
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122) # xlPasteFormats
